$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the confidentiality / as-of date text in A13.
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-10.
$ws.Range("D2").Value = 0.1076603460213693
$ws.Range("E2").Value = -0.05150022391401698

$ws.Range("D3").Value = 0.1059268734094846
$ws.Range("E3").Value = -0.0247515469716858

$ws.Range("D4").Value = 0.1160903773100152
$ws.Range("E4").Value = -0.003917646078186232

$ws.Range("D5").Value = 0.136485121431921
$ws.Range("E5").Value = -0.005920663114268865

$ws.Range("D6").Value = 0.1333811342381396
$ws.Range("E6").Value = -0.002049180327868716

$ws.Range("D7").Value = 0.1444138360402323
$ws.Range("E7").Value = -0.009084625259208168

$ws.Range("D8").Value = 0.1270556823216939
$ws.Range("E8").Value = -0.008912108174554345

$ws.Range("D9").Value = 0.1289866292271441
$ws.Range("E9").Value = -0.01581641408117052

$ws.Range("E10").Value = -0.01418697684790537
